{"js": "// Replace the 100 arithmetic answers in the single 20x5 table, in document\n// (row-major) order. The document has duplicate equation strings (e.g.\n// \"2+54=56\" and \"56+10=66\" each occur twice but map to *different* new\n// values depending on position), so replacement must be strictly\n// positional -- NOT a global find/replace keyed on old text.\nconst REPLACEMENTS = [[\"46+50=96\", \"74+11=85\"], [\"68+21=89\", \"95-67=28\"], [\"12+60=72\", \"80-43=37\"], [\"36+0=36\", \"98-70=28\"], [\"4-3=1\", \"44+25=69\"], [\"39+54=93\", \"60-33=27\"], [\"6+75=81\", \"8+16=24\"], [\"93-87=6\", \"68-61=7\"], [\"48+49=97\", \"51-22=29\"], [\"93-7=86\", \"42-11=31\"], [\"56+5=61\", \"26+52=78\"], [\"19+60=79\", \"4+56=60\"], [\"55-1=54\", \"12+29=41\"], [\"80+12=92\", \"46+16=62\"], [\"26+7=33\", \"62-13=49\"], [\"91-59=32\", \"74-63=11\"], [\"56-30=26\", \"9-8=1\"], [\"98-34=64\", \"66-50=16\"], [\"70+12=82\", \"80-46=34\"], [\"6+19=25\", \"3+62=65\"], [\"74-31=43\", \"68+23=91\"], [\"91-57=34\", \"85-18=67\"], [\"30+0=30\", \"4+74=78\"], [\"63-34=29\", \"3+43=46\"], [\"56+36=92\", \"66+11=77\"], [\"66-31=35\", \"68-52=16\"], [\"44-30=14\", \"52-45=7\"], [\"6+39=45\", \"5+92=97\"], [\"41-1=40\", \"75-20=55\"], [\"77-50=27\", \"92-83=9\"], [\"9+22=31\", \"79-21=58\"], [\"92-16=76\", \"42+23=65\"], [\"7+36=43\", \"76+0=76\"], [\"65+15=80\", \"80-68=12\"], [\"92-55=37\", \"52+5=57\"], [\"57-8=49\", \"60+37=97\"], [\"8+49=57\", \"68+5=73\"], [\"9+83=92\", \"8+51=59\"], [\"65-18=47\", \"82-52=30\"], [\"58+36=94\", \"44+2=46\"], [\"71+20=91\", \"32+44=76\"], [\"38+27=65\", \"33+26=59\"], [\"22+6=28\", \"99-24=75\"], [\"90-8=82\", \"63-17=46\"], [\"53-41=12\", \"81-43=38\"], [\"86-20=66\", \"27+59=86\"], [\"82+7=89\", \"6+25=31\"], [\"6+29=35\", \"89-59=30\"], [\"5+81=86\", \"95-55=40\"], [\"16+30=46\", \"76-0=76\"], [\"23+62=85\", \"44-17=27\"], [\"81-9=72\", \"93-0=93\"], [\"52-40=12\", \"56-41=15\"], [\"68-64=4\", \"83-33=50\"], [\"2+54=56\", \"90-26=64\"], [\"90-61=29\", \"78-76=2\"], [\"69-19=50\", \"75+13=88\"], [\"16+80=96\", \"15+59=74\"], [\"73-47=26\", \"29-9=20\"], [\"24-21=3\", \"26+39=65\"], [\"19+68=87\", \"68-6=62\"], [\"74+11=85\", \"4+9=13\"], [\"0+20=20\", \"3+30=33\"], [\"56+10=66\", \"6-4=2\"], [\"1+74=75\", \"36+48=84\"], [\"49+27=76\", \"38+32=70\"], [\"19-11=8\", \"43-42=1\"], [\"30+53=83\", \"1+6=7\"], [\"14+44=58\", \"42-2=40\"], [\"64-15=49\", \"30+42=72\"], [\"73-68=5\", \"90-77=13\"], [\"74-28=46\", \"84-47=37\"], [\"10+85=95\", \"68+0=68\"], [\"55+13=68\", \"77-11=66\"], [\"77-62=15\", \"38+9=47\"], [\"54+34=88\", \"94-65=29\"], [\"2+54=56\", \"72+0=72\"], [\"46+17=63\", \"70-11=59\"], [\"29-22=7\", \"62-6=56\"], [\"56+10=66\", \"26+20=46\"], [\"43+37=80\", \"11+60=71\"], [\"7+31=38\", \"28-5=23\"], [\"92-90=2\", \"75-44=31\"], [\"3+12=15\", \"38-4=34\"], [\"96-39=57\", \"67-20=47\"], [\"40-31=9\", \"31+6=37\"], [\"62-23=39\", \"3+9=12\"], [\"79+15=94\", \"28+47=75\"], [\"8+17=25\", \"43+5=48\"], [\"11+19=30\", \"86-67=19\"], [\"35+6=41\", \"27-11=16\"], [\"99-6=93\", \"40+26=66\"], [\"33+9=42\", \"9+23=32\"], [\"62+8=70\", \"94-89=5\"], [\"65-35=30\", \"44+55=99\"], [\"96-51=45\", \"26+6=32\"], [\"12+30=42\", \"27+48=75\"], [\"69+4=73\", \"26+64=90\"], [\"25-14=11\", \"22+5=27\"], [\"1+50=51\", \"93+1=94\"]];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst oldValues = table.values;\nconst rows = oldValues.length;\nconst cols = oldValues[0].length;\n\nconst newValues = [];\nlet k = 0;\nfor (let r = 0; r < rows; r++) {\n  const newRow = [];\n  for (let c = 0; c < cols; c++) {\n    const cellText = oldValues[r][c];\n    if (k < REPLACEMENTS.length && cellText === REPLACEMENTS[k][0]) {\n      newRow.push(REPLACEMENTS[k][1]);\n    } else {\n      // Fallback: keep the original text untouched if something doesn't\n      // line up the way we expect, rather than corrupting unrelated cells.\n      newRow.push(cellText);\n    }\n    k++;\n  }\n  newValues.push(newRow);\n}\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Replace the 100 arithmetic answers in the single 20x5 table, in document\n# (row-major) order. The document has duplicate equation strings (e.g.\n# \"2+54=56\" and \"56+10=66\" each occur twice but map to *different* new\n# values depending on position), so replacement must be strictly\n# positional -- NOT a global find/replace keyed on old text.\n$REPLACEMENTS = @(\n    @(\"46+50=96\", \"74+11=85\"),\n    @(\"68+21=89\", \"95-67=28\"),\n    @(\"12+60=72\", \"80-43=37\"),\n    @(\"36+0=36\", \"98-70=28\"),\n    @(\"4-3=1\", \"44+25=69\"),\n    @(\"39+54=93\", \"60-33=27\"),\n    @(\"6+75=81\", \"8+16=24\"),\n    @(\"93-87=6\", \"68-61=7\"),\n    @(\"48+49=97\", \"51-22=29\"),\n    @(\"93-7=86\", \"42-11=31\"),\n    @(\"56+5=61\", \"26+52=78\"),\n    @(\"19+60=79\", \"4+56=60\"),\n    @(\"55-1=54\", \"12+29=41\"),\n    @(\"80+12=92\", \"46+16=62\"),\n    @(\"26+7=33\", \"62-13=49\"),\n    @(\"91-59=32\", \"74-63=11\"),\n    @(\"56-30=26\", \"9-8=1\"),\n    @(\"98-34=64\", \"66-50=16\"),\n    @(\"70+12=82\", \"80-46=34\"),\n    @(\"6+19=25\", \"3+62=65\"),\n    @(\"74-31=43\", \"68+23=91\"),\n    @(\"91-57=34\", \"85-18=67\"),\n    @(\"30+0=30\", \"4+74=78\"),\n    @(\"63-34=29\", \"3+43=46\"),\n    @(\"56+36=92\", \"66+11=77\"),\n    @(\"66-31=35\", \"68-52=16\"),\n    @(\"44-30=14\", \"52-45=7\"),\n    @(\"6+39=45\", \"5+92=97\"),\n    @(\"41-1=40\", \"75-20=55\"),\n    @(\"77-50=27\", \"92-83=9\"),\n    @(\"9+22=31\", \"79-21=58\"),\n    @(\"92-16=76\", \"42+23=65\"),\n    @(\"7+36=43\", \"76+0=76\"),\n    @(\"65+15=80\", \"80-68=12\"),\n    @(\"92-55=37\", \"52+5=57\"),\n    @(\"57-8=49\", \"60+37=97\"),\n    @(\"8+49=57\", \"68+5=73\"),\n    @(\"9+83=92\", \"8+51=59\"),\n    @(\"65-18=47\", \"82-52=30\"),\n    @(\"58+36=94\", \"44+2=46\"),\n    @(\"71+20=91\", \"32+44=76\"),\n    @(\"38+27=65\", \"33+26=59\"),\n    @(\"22+6=28\", \"99-24=75\"),\n    @(\"90-8=82\", \"63-17=46\"),\n    @(\"53-41=12\", \"81-43=38\"),\n    @(\"86-20=66\", \"27+59=86\"),\n    @(\"82+7=89\", \"6+25=31\"),\n    @(\"6+29=35\", \"89-59=30\"),\n    @(\"5+81=86\", \"95-55=40\"),\n    @(\"16+30=46\", \"76-0=76\"),\n    @(\"23+62=85\", \"44-17=27\"),\n    @(\"81-9=72\", \"93-0=93\"),\n    @(\"52-40=12\", \"56-41=15\"),\n    @(\"68-64=4\", \"83-33=50\"),\n    @(\"2+54=56\", \"90-26=64\"),\n    @(\"90-61=29\", \"78-76=2\"),\n    @(\"69-19=50\", \"75+13=88\"),\n    @(\"16+80=96\", \"15+59=74\"),\n    @(\"73-47=26\", \"29-9=20\"),\n    @(\"24-21=3\", \"26+39=65\"),\n    @(\"19+68=87\", \"68-6=62\"),\n    @(\"74+11=85\", \"4+9=13\"),\n    @(\"0+20=20\", \"3+30=33\"),\n    @(\"56+10=66\", \"6-4=2\"),\n    @(\"1+74=75\", \"36+48=84\"),\n    @(\"49+27=76\", \"38+32=70\"),\n    @(\"19-11=8\", \"43-42=1\"),\n    @(\"30+53=83\", \"1+6=7\"),\n    @(\"14+44=58\", \"42-2=40\"),\n    @(\"64-15=49\", \"30+42=72\"),\n    @(\"73-68=5\", \"90-77=13\"),\n    @(\"74-28=46\", \"84-47=37\"),\n    @(\"10+85=95\", \"68+0=68\"),\n    @(\"55+13=68\", \"77-11=66\"),\n    @(\"77-62=15\", \"38+9=47\"),\n    @(\"54+34=88\", \"94-65=29\"),\n    @(\"2+54=56\", \"72+0=72\"),\n    @(\"46+17=63\", \"70-11=59\"),\n    @(\"29-22=7\", \"62-6=56\"),\n    @(\"56+10=66\", \"26+20=46\"),\n    @(\"43+37=80\", \"11+60=71\"),\n    @(\"7+31=38\", \"28-5=23\"),\n    @(\"92-90=2\", \"75-44=31\"),\n    @(\"3+12=15\", \"38-4=34\"),\n    @(\"96-39=57\", \"67-20=47\"),\n    @(\"40-31=9\", \"31+6=37\"),\n    @(\"62-23=39\", \"3+9=12\"),\n    @(\"79+15=94\", \"28+47=75\"),\n    @(\"8+17=25\", \"43+5=48\"),\n    @(\"11+19=30\", \"86-67=19\"),\n    @(\"35+6=41\", \"27-11=16\"),\n    @(\"99-6=93\", \"40+26=66\"),\n    @(\"33+9=42\", \"9+23=32\"),\n    @(\"62+8=70\", \"94-89=5\"),\n    @(\"65-35=30\", \"44+55=99\"),\n    @(\"96-51=45\", \"26+6=32\"),\n    @(\"12+30=42\", \"27+48=75\"),\n    @(\"69+4=73\", \"26+64=90\"),\n    @(\"25-14=11\", \"22+5=27\"),\n    @(\"1+50=51\", \"93+1=94\")\n)\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n$rowCount = $tbl.Rows.Count\n$colCount = $tbl.Columns.Count\n\n$k = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $tbl.Cell($r, $c)\n        $pair = $REPLACEMENTS[$k]\n        $oldExpected = $pair[0]\n        $newValue = $pair[1]\n\n        $cellRange = $cell.Range\n        # Trim the trailing end-of-cell marker(s) (CR + cell-mark) before\n        # comparing so we match on the visible text only.\n        $currentText = $cellRange.Text.TrimEnd([char]13, [char]7)\n\n        if ($currentText -eq $oldExpected) {\n            $cellRange.Text = $newValue\n        }\n        # else: leave the cell untouched rather than risk corrupting an\n        # unrelated cell if the table layout doesn't match expectations.\n\n        $k++\n    }\n}\n"}
